# Update "想去人数" (interested-people count) figures in both the
# "展览" sheet and the "全部类型" sheet (which mirrors the same data),
# matching the latest scrape output.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F3").Value = 6291
    $ws.Range("F8").Value = 1433
    $ws.Range("F11").Value = 256
}
